$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = 21.63579072432197
$ws.Cells.Item(2,2).Value = 17.17264978995248
$ws.Cells.Item(2,3).Value = 25.96870934173277
$ws.Cells.Item(3,1).Value = 20.70305128798127
$ws.Cells.Item(3,2).Value = 14.91734937250084
$ws.Cells.Item(3,3).Value = 26.79243267464477
$ws.Cells.Item(4,1).Value = 20.47224903958168
$ws.Cells.Item(4,2).Value = 15.03354342944503
$ws.Cells.Item(4,3).Value = 26.53866393881437
$ws.Cells.Item(5,1).Value = 25.07756018209518
$ws.Cells.Item(5,2).Value = 17.23966322089311
$ws.Cells.Item(5,3).Value = 33.38535189571672
$ws.Cells.Item(6,1).Value = 26.30259005123169
$ws.Cells.Item(6,2).Value = 20.99772675760862
$ws.Cells.Item(6,3).Value = 32.23333504973511
$ws.Cells.Item(7,1).Value = 19.2395070158974
$ws.Cells.Item(7,2).Value = 15.56094316134581
$ws.Cells.Item(7,3).Value = 22.97759726621279
$ws.Cells.Item(8,1).Value = 9.758607296095168
$ws.Cells.Item(8,2).Value = 4.572337952333323
$ws.Cells.Item(8,3).Value = 15.48420273623262
$ws.Cells.Item(9,1).Value = 24.52047265564187
$ws.Cells.Item(9,2).Value = 19.54869120689143
$ws.Cells.Item(9,3).Value = 29.62870346239988
$ws.Cells.Item(10,1).Value = 10.04562654708925
$ws.Cells.Item(10,2).Value = 4.874639802632513
$ws.Cells.Item(10,3).Value = 16.566092765304
$ws.Cells.Item(11,1).Value = 18.1973433473162
$ws.Cells.Item(11,2).Value = 13.25037058664995
$ws.Cells.Item(11,3).Value = 23.06925489301352
$ws.Cells.Item(12,1).Value = 27.02967906618417
$ws.Cells.Item(12,2).Value = 21.00944182041341
$ws.Cells.Item(12,3).Value = 33.03499597986732
$ws.Cells.Item(13,1).Value = 13.08656409466251
$ws.Cells.Item(13,2).Value = 6.886715415386097
$ws.Cells.Item(13,3).Value = 20.64969630535904
$ws.Cells.Item(14,1).Value = 14.92332544723742
$ws.Cells.Item(14,2).Value = 9.756495689784401
$ws.Cells.Item(14,3).Value = 20.63271830802748
$ws.Cells.Item(15,1).Value = 29.10299263033501
$ws.Cells.Item(15,2).Value = 21.87953785492587
$ws.Cells.Item(15,3).Value = 37.90948306225241
$ws.Cells.Item(16,1).Value = 26.34848384992639
$ws.Cells.Item(16,2).Value = 20.78918657120634
$ws.Cells.Item(16,3).Value = 32.40085757335024
$ws.Cells.Item(17,1).Value = 26.92823095490016
$ws.Cells.Item(17,2).Value = 20.91611222597446
$ws.Cells.Item(17,3).Value = 32.54737131527081
$ws.Cells.Item(18,1).Value = 20.24562036735332
$ws.Cells.Item(18,2).Value = 13.78319719821131
$ws.Cells.Item(18,3).Value = 26.62160544793439
$ws.Cells.Item(19,1).Value = 18.80334786968479
$ws.Cells.Item(19,2).Value = 13.18515640161169
$ws.Cells.Item(19,3).Value = 24.82297478523168
$ws.Cells.Item(20,1).Value = 29.58992741740598
$ws.Cells.Item(20,2).Value = 22.90467894907736
$ws.Cells.Item(20,3).Value = 35.7568351426146
$ws.Cells.Item(21,1).Value = 29.18128709902766
$ws.Cells.Item(21,2).Value = 21.85392684109442
$ws.Cells.Item(21,3).Value = 37.18851797693458
$ws.Cells.Item(22,1).Value = 10.56451300613945
$ws.Cells.Item(22,2).Value = 4.702303892045861
$ws.Cells.Item(22,3).Value = 17.57865338323183
$ws.Cells.Item(23,1).Value = 31.22799584798681
$ws.Cells.Item(23,2).Value = 22.93716838578001
$ws.Cells.Item(23,3).Value = 40.44049242433699
$ws.Cells.Item(24,1).Value = 26.34424863221215
$ws.Cells.Item(24,2).Value = 20.71442796138463
$ws.Cells.Item(24,3).Value = 32.62988860297703
$ws.Cells.Item(25,1).Value = 22.60926801062947
$ws.Cells.Item(25,2).Value = 17.4090672753626
$ws.Cells.Item(25,3).Value = 27.30860046120073
$ws.Cells.Item(26,1).Value = 29.04473888346227
$ws.Cells.Item(26,2).Value = 22.1457075098416
$ws.Cells.Item(26,3).Value = 36.32034412974819
$ws.Cells.Item(27,1).Value = 12.9727590119992
$ws.Cells.Item(27,2).Value = 6.874612378818136
$ws.Cells.Item(27,3).Value = 21.27482847456812
$ws.Cells.Item(28,1).Value = 20.42249345384324
$ws.Cells.Item(28,2).Value = 16.32999567783657
$ws.Cells.Item(28,3).Value = 24.31228631483254
$ws.Cells.Item(29,1).Value = 32.28469923388152
$ws.Cells.Item(29,2).Value = 22.78672834044227
$ws.Cells.Item(29,3).Value = 41.38715332103755
$ws.Cells.Item(30,1).Value = 25.5643224126752
$ws.Cells.Item(30,2).Value = 20.92225697363093
$ws.Cells.Item(30,3).Value = 31.16730472899639
$ws.Cells.Item(31,1).Value = 14.22364776780775
$ws.Cells.Item(31,2).Value = 9.069700886549185
$ws.Cells.Item(31,3).Value = 20.27273505158202
$ws.Cells.Item(32,1).Value = 11.9912383611151
$ws.Cells.Item(32,2).Value = 6.724831146663179
$ws.Cells.Item(32,3).Value = 17.54209010921339
$ws.Cells.Item(33,1).Value = 14.90108117739871
$ws.Cells.Item(33,2).Value = 9.784527622641269
$ws.Cells.Item(33,3).Value = 20.59819368050134
$ws.Cells.Item(34,1).Value = 23.31559126682208
$ws.Cells.Item(34,2).Value = 18.26855820413786
$ws.Cells.Item(34,3).Value = 29.0365660872575
$ws.Cells.Item(35,1).Value = 28.03538873014233
$ws.Cells.Item(35,2).Value = 22.70945777714157
$ws.Cells.Item(35,3).Value = 33.87922883820316
$ws.Cells.Item(36,1).Value = 19.26421954552482
$ws.Cells.Item(36,2).Value = 15.38008130853053
$ws.Cells.Item(36,3).Value = 22.9361467076763
$ws.Cells.Item(37,1).Value = 31.1247967342051
$ws.Cells.Item(37,2).Value = 23.88125828394762
$ws.Cells.Item(37,3).Value = 39.33016966495096
$ws.Cells.Item(38,1).Value = 15.30417375195411
$ws.Cells.Item(38,2).Value = 9.798207505718823
$ws.Cells.Item(38,3).Value = 22.2977138269052
$ws.Cells.Item(39,1).Value = 27.69083953003661
$ws.Cells.Item(39,2).Value = 22.55457272277003
$ws.Cells.Item(39,3).Value = 33.4192396935374
$ws.Cells.Item(40,1).Value = 25.17647674168717
$ws.Cells.Item(40,2).Value = 18.40364953345121
$ws.Cells.Item(40,3).Value = 32.65529260606988
$ws.Cells.Item(41,1).Value = 18.19517717317779
$ws.Cells.Item(41,2).Value = 14.77968062922
$ws.Cells.Item(41,3).Value = 21.85691225496999
$ws.Cells.Item(42,1).Value = 18.18933494570941
$ws.Cells.Item(42,2).Value = 13.32134815846064
$ws.Cells.Item(42,3).Value = 23.36069320009562
$ws.Cells.Item(43,1).Value = 25.5643224126752
$ws.Cells.Item(43,2).Value = 20.92225697363093
$ws.Cells.Item(43,3).Value = 31.16730472899639
$ws.Cells.Item(44,1).Value = 16.01514517925195
$ws.Cells.Item(44,2).Value = 12.32129761222182
$ws.Cells.Item(44,3).Value = 19.39032600296134
$ws.Cells.Item(45,1).Value = 19.07194974413987
$ws.Cells.Item(45,2).Value = 13.57937445145544
$ws.Cells.Item(45,3).Value = 25.16900287873549
$ws.Cells.Item(46,1).Value = 13.71176446213956
$ws.Cells.Item(46,2).Value = 9.739883332439929
$ws.Cells.Item(46,3).Value = 18.53428137365061
$ws.Cells.Item(47,1).Value = 10.94437837105579
$ws.Cells.Item(47,2).Value = 5.044763133006941
$ws.Cells.Item(47,3).Value = 18.53747801741599
$ws.Cells.Item(48,1).Value = 27.04747547978573
$ws.Cells.Item(48,2).Value = 21.04455052559879
$ws.Cells.Item(48,3).Value = 33.05795484521368
$ws.Cells.Item(49,1).Value = 22.99826433434714
$ws.Cells.Item(49,2).Value = 17.3326328690722
$ws.Cells.Item(49,3).Value = 29.68548834607849
$ws.Cells.Item(50,1).Value = 31.18020859211034
$ws.Cells.Item(50,2).Value = 23.95563841205932
$ws.Cells.Item(50,3).Value = 39.43473436036884
$ws.Cells.Item(51,1).Value = 20.38992277549125
$ws.Cells.Item(51,2).Value = 14.64651293459907
$ws.Cells.Item(51,3).Value = 26.36253094699647
$ws.Cells.Item(52,1).Value = 13.56976321778885
$ws.Cells.Item(52,2).Value = 9.392428169503807
$ws.Cells.Item(52,3).Value = 17.83816703112846
$ws.Cells.Item(53,1).Value = 27.79721683453665
$ws.Cells.Item(53,2).Value = 22.77454519396848
$ws.Cells.Item(53,3).Value = 33.51112941012799
$ws.Cells.Item(54,1).Value = 29.00812814866874
$ws.Cells.Item(54,2).Value = 22.04964256772849
$ws.Cells.Item(54,3).Value = 37.90651784225873
$ws.Cells.Item(55,1).Value = 28.96679446127066
$ws.Cells.Item(55,2).Value = 22.00320938952697
$ws.Cells.Item(55,3).Value = 36.73060989850141
$ws.Cells.Item(56,1).Value = 14.50262704229279
$ws.Cells.Item(56,2).Value = 9.644026968491946
$ws.Cells.Item(56,3).Value = 19.69007915938051
$ws.Cells.Item(57,1).Value = 11.95629682071792
$ws.Cells.Item(57,2).Value = 6.620721681949081
$ws.Cells.Item(57,3).Value = 17.73387617593886
$ws.Cells.Item(58,1).Value = 31.02747577862743
$ws.Cells.Item(58,2).Value = 23.60316328057241
$ws.Cells.Item(58,3).Value = 39.39524836420587
$ws.Cells.Item(59,1).Value = 17.06546016391139
$ws.Cells.Item(59,2).Value = 13.05427526007521
$ws.Cells.Item(59,3).Value = 21.3638585302112
$ws.Cells.Item(60,1).Value = 31.61753489926865
$ws.Cells.Item(60,2).Value = 22.69417792274855
$ws.Cells.Item(60,3).Value = 40.53601238892934
$ws.Cells.Item(61,1).Value = 15.11296293230981
$ws.Cells.Item(61,2).Value = 9.806438085000313
$ws.Cells.Item(61,3).Value = 21.25231941398624
$ws.Cells.Item(62,1).Value = 18.20385468471127
$ws.Cells.Item(62,2).Value = 13.36583750647366
$ws.Cells.Item(62,3).Value = 22.93417247425351
$ws.Cells.Item(63,1).Value = 15.31517170986328
$ws.Cells.Item(63,2).Value = 8.96982502303778
$ws.Cells.Item(63,3).Value = 22.76488566537851
$ws.Cells.Item(64,1).Value = 31.41445391729707
$ws.Cells.Item(64,2).Value = 22.1938641061934
$ws.Cells.Item(64,3).Value = 40.32325371330696
$ws.Cells.Item(65,1).Value = 11.44779292206309
$ws.Cells.Item(65,2).Value = 6.726251378379641
$ws.Cells.Item(65,3).Value = 16.50833018745817
$ws.Cells.Item(66,1).Value = 22.97041411049093
$ws.Cells.Item(66,2).Value = 17.66206578859173
$ws.Cells.Item(66,3).Value = 28.49334438955825
$ws.Cells.Item(67,1).Value = 22.07311381888328
$ws.Cells.Item(67,2).Value = 17.5288409348339
$ws.Cells.Item(67,3).Value = 27.03361231410391
$ws.Cells.Item(68,1).Value = 26.82506545755995
$ws.Cells.Item(68,2).Value = 19.8985058540083
$ws.Cells.Item(68,3).Value = 35.48769371613449
$ws.Cells.Item(69,1).Value = 11.17996538183067
$ws.Cells.Item(69,2).Value = 5.713358385685818
$ws.Cells.Item(69,3).Value = 18.77789098213669
$ws.Cells.Item(70,1).Value = 25.05754487370545
$ws.Cells.Item(70,2).Value = 18.47472214930172
$ws.Cells.Item(70,3).Value = 32.08634126392594
$ws.Cells.Item(71,1).Value = 19.77501335341309
$ws.Cells.Item(71,2).Value = 16.49828977680108
$ws.Cells.Item(71,3).Value = 23.8552482296575
$ws.Cells.Item(72,1).Value = 23.02678608992188
$ws.Cells.Item(72,2).Value = 18.97100883837469
$ws.Cells.Item(72,3).Value = 27.78141462492555
